$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the scraped price/volume (and, for rows 32-33, the row-swapped
# coin identity) updates cell by cell, in sheet order.
#
# Cells whose new text reads as a plain decimal number (e.g. "1.013",
# "0.09932") would otherwise be auto-converted to a numeric value by
# Excel's normal type inference on Range.Value assignment, which would
# lose the original text formatting (leading zeros / fixed decimal counts)
# that the source file stores as plain text. For those, force the cell to
# Text format before writing, then clear the now-unneeded explicit format
# so the cell's style index matches the untouched cells around it.

$ws.Range("D2").Value = '29.574.58'
$ws.Range("E2").Value = '  -2.34%  '
$ws.Range("D3").Value = '2.000.43'
$ws.Range("E3").Value = '  -4.27%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +1.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.37'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.05%  '
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5001'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -4.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4226'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.92'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08988'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.70%  '
$ws.Range("E11").Value = '  -4.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.33'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -5.79%  '
$ws.Range("D13").Value = '2.035.35'
$ws.Range("E13").Value = '  -2.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.059'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -6.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.474'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -6.13%  '
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.03'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -6.88%  '
$ws.Range("E18").Value = '  -3.81%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06676'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.70'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -6.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.011'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.953'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -6.26%  '
$ws.Range("D23").Value = '29.591.31'
$ws.Range("E23").Value = '  -2.23%  '
$ws.Range("E24").Value = '  -4.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.301'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.51'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.30%  '
$ws.Range("E27").Value = '  -5.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.420'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.300'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -8.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.15'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.53%  '
$ws.Range("E31").Value = '  -6.76%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09932'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.15%  '
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.574'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.08%  '
$ws.Range("E34").Value = '  -6.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.800'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02467'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.306'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -8.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.306'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06347'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6558'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -6.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.67'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -6.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2048'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.010'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6342'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.50'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.198'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.307'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.503'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000342'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06984'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.33%  '
$ws.Range("E51").Value = '  -7.49%  '
